$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap currency labels for rows 4 and 5: BTC now appears before ATOM
$ws.Range("A4").Value = "BTC"
$ws.Range("A5").Value = "ATOM"

# Row 3 (USDT) - refreshed balances
$ws.Range("C3").Value = 863.0471265
$ws.Range("D3").Value = 90.62048405
$ws.Range("E3").Value = 772.42664245
$ws.Range("G3").Value = 863.05

# Row 4 (now BTC) - values from former BTC row, with refreshed price/dollar_value
$ws.Range("C4").Value = 0.01540148
$ws.Range("D4").Value = 0.01540148
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 33509
$ws.Range("G4").Value = 516.09

# Row 5 (now ATOM) - refreshed balances and price/dollar_value
$ws.Range("C5").Value = 24.0438
$ws.Range("D5").Value = 0.2372
$ws.Range("E5").Value = 23.8066
$ws.Range("F5").Value = 12.7
$ws.Range("G5").Value = 305.36

# Row 6 (ALGO) - refreshed price only
$ws.Range("F6").Value = 0.8508

# Row 7 (ETH) - refreshed price only
$ws.Range("F7").Value = 2211.73
